# Stats Table - QOL improvements and bug fixes
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "Results" sheet: drop the finished match rows (the "cris" matches and the
#    lost "test" match), keeping only the still-pending "test / NA" row (was
#    row 5) which shifts up to become row 2.
# ---------------------------------------------------------------------------
$results = $wb.Worksheets.Item("Results")

$results.Rows.Item(4).Delete()
$results.Rows.Item(3).Delete()
$results.Rows.Item(2).Delete()

# ---------------------------------------------------------------------------
# 2. "h2h" sheet: the head-to-head record against "test" resets to 0-0 and
#    the now-duplicate "test" summary row is removed.
# ---------------------------------------------------------------------------
$h2h = $wb.Worksheets.Item("h2h")

# Write "0" as text (matching the existing text-stored "1"/"0" counters)
# without leaving a quote-prefix number format behind: type it with a
# leading apostrophe, then re-pull the plain/general format from a
# neighbouring untouched cell.
$h2h.Range("A2").Formula = "'0"
$h2h.Range("B2").Copy()
$h2h.Range("A2").PasteSpecial(-4122) | Out-Null

$h2h.Range("C2").Formula = "'0"
$h2h.Range("D2").Copy()
$h2h.Range("C2").PasteSpecial(-4122) | Out-Null

$h2h.Range("E2").Value = "test"

# Drop the old "test" row (row 3); the already-blank trailing row shifts up
# to become the new blank row 3.
$h2h.Rows.Item(3).Delete()

# ---------------------------------------------------------------------------
# 3. Add a new "stats" sheet after "h2h" with the all-time record breakdown.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$stats = $wb.Worksheets.Add($null, $lastSheet)
$stats.Name = "stats"

$stats.Range("B1").Value = "All Time"

$stats.Range("A2").Value = "Overall"
$stats.Range("B2").Value = "0-0"

$stats.Range("A3").Value = "Clay"
$stats.Range("B3").Value = "0-0"

$stats.Range("A4").Value = "Hard"
$stats.Range("B4").Value = "0-0"

$stats.Range("A5").Value = "Tartan"
$stats.Range("B5").Value = "0-0"

# Make the new "stats" sheet the active tab.
$stats.Activate()
